$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells H1 and I1, matching the style of existing header cells
$ws.Range("G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)

$ws.Range("H1").Value = "Situação"
$ws.Range("I1").Value = "Status"

# Update row 2 values
$ws.Range("A2").Value = 23003
$ws.Range("B2").Value = "RCR-BA"
$ws.Range("C2").Value = "16 de fev. de 2025, 17:09:29"
$ws.Range("F2").Value = "https://moovsec-videos-prod.s3.sa-east-1.amazonaws.com/video_evidence_67b2461209aca087c3be2420.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Date=20250217T063634Z&X-Amz-SignedHeaders=host&X-Amz-Expires=604800&X-Amz-Credential=AKIA3QTD5B6Z3DVDF6WA%2F20250217%2Fsa-east-1%2Fs3%2Faws4_request&X-Amz-Signature=1f449535706cc21f830890e984cf99ef4c6be9ee5793ac8c5e0be28c0f91a0ba"
$ws.Range("G2").Value = "https://www.google.com/maps?q=-12.44133900000000000,-37.92349799999999500"
$ws.Range("H2").Value = "Falso"
$ws.Range("I2").Value = "Validado"

# Remove rows 3 and 4 (no longer present in the data)
$ws.Rows("3:4").Delete()
